$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> @(newPriceOrNull, newVolumePct)
$updates = @{
    2 = @("26.998.90", "  +0.57%  ")
    3 = @("1.640.24", "  -0.01%  ")
    4 = @($null, "  -0.53%  ")
    5 = @("218.07", "  +0.12%  ")
    6 = @($null, "  +1.97%  ")
    7 = @($null, "  -0.51%  ")
    8 = @($null, "  +1.62%  ")
    9 = @($null, "  +0.38%  ")
    10 = @("20.00", "  +3.90%  ")
    11 = @("0.0846", "  +0.17%  ")
    12 = @("1.869.29", "  -0.02%  ")
    13 = @("1.631.55", "  -0.49%  ")
    14 = @($null, "  -0.86%  ")
    15 = @("0.534", "  +1.37%  ")
    16 = @("67.15", "  +2.95%  ")
    17 = @("26.977.72", "  +0.53%  ")
    18 = @($null, "  +0.30%  ")
    19 = @("219.93", "  +2.16%  ")
    20 = @($null, "  -0.55%  ")
    21 = @($null, "  +3.15%  ")
    22 = @($null, "  +1.30%  ")
    23 = @($null, "  +1.62%  ")
    24 = @("9.19", "  -0.10%  ")
    25 = @("147.41", "  +0.08%  ")
    26 = @($null, "  -0.44%  ")
    27 = @("7.35", "  +2.18%  ")
    28 = @($null, "  +1.04%  ")
    29 = @("15.77", "  +0.25%  ")
    30 = @("0.0505", "  -0.43%  ")
    31 = @($null, "  -0.32%  ")
    32 = @($null, "  -0.70%  ")
    33 = @($null, "  +0.56%  ")
    34 = @($null, "  +1.29%  ")
    35 = @("1.270.04", "  -0.17%  ")
    36 = @($null, "  -0.23%  ")
    37 = @($null, "  +2.63%  ")
    38 = @("0.544", "  +2.70%  ")
    39 = @($null, "  +2.65%  ")
    40 = @($null, "  -0.47%  ")
    41 = @("0.807", "  +0.51%  ")
    42 = @("5.36", "  +0.80%  ")
    43 = @("1.780.38", "  +0.03%  ")
    44 = @($null, "  +3.59%  ")
    45 = @("62.25", "  +2.13%  ")
    46 = @("92.55", "  +0.15%  ")
    47 = @($null, "  +1.72%  ")
    48 = @("0.0₆0105", "  +0.48%  ")
    49 = @($null, "  -0.53%  ")
    50 = @("7.69", "  +1.71%  ")
    51 = @("0.0966", "  -0.08%  ")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $newPrice = $vals[0]
    $newVolume = $vals[1]
    if ($null -ne $newPrice) {
        $priceCell = $ws.Cells.Item($row, 4)
        if ($newPrice -match "^-?[0-9]+(\.[0-9]+)?$") {
            # Looks like a plain number to Excel auto-detection; force text so it
            # keeps trailing zeros / fixed decimal places exactly as scraped.
            $priceCell.NumberFormat = "@"
        }
        $priceCell.Value = $newPrice
    }
    $ws.Cells.Item($row, 5).Value = $newVolume
}
